$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 95
$ws.Range("J9").Value = 75
$ws.Range("L9").Value = 75
$ws.Range("N9").Value = -413
$ws.Range("H129").Value = 1334.8889
$ws.Range("I129").Value = 572
$ws.Range("J129").Value = 1945.2
$ws.Range("K129").Value = 1716
$ws.Range("L129").Value = 5835.6
$ws.Range("M129").Value = 3284
$ws.Range("N129").Value = -15835.6
$ws.Range("H132").Value = 238379.64
$ws.Range("I132").Value = 259451.81
$ws.Range("J132").Value = 40301.2
$ws.Range("K132").Value = 778355.4299999999
$ws.Range("L132").Value = 120903.6
$ws.Range("M132").Value = -775825.4299999999
$ws.Range("N132").Value = -125963.6
$ws.Range("H135").Value = 1208.0667
$ws.Range("I135").Value = 1089.0927
$ws.Range("J135").Value = 1514
$ws.Range("K135").Value = 9801.834299999999
$ws.Range("L135").Value = 13626
$ws.Range("M135").Value = -7266.834299999999
$ws.Range("N135").Value = -18696
$ws.Range("H137").Value = 16129979
$ws.Range("I137").Value = 20000698
$ws.Range("J137").Value = 1981.4166
$ws.Range("K137").Value = 60002094
$ws.Range("L137").Value = 5944.2498
$ws.Range("M137").Value = -59999544
$ws.Range("N137").Value = -11044.2498
$ws.Range("H138").Value = 8334911.5
$ws.Range("I138").Value = 970561.9399999999
$ws.Range("J138").Value = 166668420
$ws.Range("K138").Value = 2911685.82
$ws.Range("L138").Value = 500005260
$ws.Range("M138").Value = -2906545.82
$ws.Range("N138").Value = -500015540

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17053.629
$ws.Range("I32").Value = 2534.4614
$ws.Range("J32").Value = 205802.8
$ws.Range("K32").Value = 2534.4614
$ws.Range("L32").Value = 205802.8
$ws.Range("M32").Value = -2247.4614
$ws.Range("N32").Value = -206376.8
$ws.Range("H61").Value = 1551.0344
$ws.Range("I61").Value = 1234.711
$ws.Range("J61").Value = 2646
$ws.Range("K61").Value = 1234.711
$ws.Range("L61").Value = 2646
$ws.Range("M61").Value = -1022.711
$ws.Range("N61").Value = -3070
$ws.Range("H74").Value = 3193.492
$ws.Range("I74").Value = 991.0732
$ws.Range("J74").Value = 7298
$ws.Range("K74").Value = 991.0732
$ws.Range("L74").Value = 7298
$ws.Range("M74").Value = -117.0732
$ws.Range("N74").Value = -9046
$ws.Range("H77").Value = 3193.492
$ws.Range("I77").Value = 991.0732
$ws.Range("J77").Value = 7298
$ws.Range("K77").Value = 4955.366
$ws.Range("L77").Value = 36490
$ws.Range("M77").Value = -587.366
$ws.Range("N77").Value = -45226
$ws.Range("H97").Value = 20833750
$ws.Range("I97").Value = 20833750
$ws.Range("K97").Value = 20833750
$ws.Range("M97").Value = -20833254
$ws.Range("H110").Value = 803.58826
$ws.Range("I110").Value = 729
$ws.Range("J110").Value = 982.6
$ws.Range("K110").Value = 729
$ws.Range("L110").Value = 982.6
$ws.Range("M110").Value = 1316
$ws.Range("N110").Value = -5072.6
$ws.Range("H122").Value = 2727.9285
$ws.Range("I122").Value = 2776.1
$ws.Range("J122").Value = 2607.5
$ws.Range("K122").Value = 8328.299999999999
$ws.Range("L122").Value = 7822.5
$ws.Range("M122").Value = -5878.299999999999
$ws.Range("N122").Value = -12722.5
$ws.Range("H132").Value = 2104.7036
$ws.Range("I132").Value = 1526.8864
$ws.Range("J132").Value = 4647.1
$ws.Range("K132").Value = 4580.6592
$ws.Range("L132").Value = 13941.3
$ws.Range("M132").Value = -2050.6592
$ws.Range("N132").Value = -19001.3
$ws.Range("H136").Value = 1551.0344
$ws.Range("I136").Value = 1234.711
$ws.Range("J136").Value = 2646
$ws.Range("K136").Value = 3704.133
$ws.Range("L136").Value = 7938
$ws.Range("M136").Value = -1154.133
$ws.Range("N136").Value = -13038

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1280.2354
$ws.Range("I94").Value = 1063.6
$ws.Range("J94").Value = 2905
$ws.Range("K94").Value = 1063.6
$ws.Range("L94").Value = 2905
$ws.Range("M94").Value = -612.5999999999999
$ws.Range("N94").Value = -3807
$ws.Range("H134").Value = 16130848
$ws.Range("I134").Value = 19609232
$ws.Range("J134").Value = 3794.5454
$ws.Range("K134").Value = 58827696
$ws.Range("L134").Value = 11383.6362
$ws.Range("M134").Value = -58825161
$ws.Range("N134").Value = -16453.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2084.8215
$ws.Range("I31").Value = 1224.7894
$ws.Range("J31").Value = 3900.4443
$ws.Range("K31").Value = 1224.7894
$ws.Range("L31").Value = 3900.4443
$ws.Range("M31").Value = -929.7893999999999
$ws.Range("N31").Value = -4490.4443
$ws.Range("H34").Value = 2084.8215
$ws.Range("I34").Value = 1224.7894
$ws.Range("J34").Value = 3900.4443
$ws.Range("K34").Value = 1224.7894
$ws.Range("L34").Value = 3900.4443
$ws.Range("M34").Value = -1022.7894
$ws.Range("N34").Value = -4304.4443

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 52631736
$ws.Range("I2").Value = 36.416668
$ws.Range("J2").Value = 142857500
$ws.Range("K2").Value = 218.500008
$ws.Range("L2").Value = 857145000
$ws.Range("M2").Value = -105.500008
$ws.Range("N2").Value = -857145226
$ws.Range("H40").Value = 437.5
$ws.Range("J40").Value = 490
$ws.Range("L40").Value = 1960
$ws.Range("N40").Value = -2098

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 43333.332
$ws.Range("I35").Value = 40000
$ws.Range("J35").Value = 50000
$ws.Range("K35").Value = 40000
$ws.Range("L35").Value = 50000
$ws.Range("M35").Value = -39702
$ws.Range("N35").Value = -50596
$ws.Range("H40").Value = 8000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 8000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -8302
$ws.Range("H132").Value = 2617.5588
$ws.Range("I132").Value = 2130.0188
$ws.Range("J132").Value = 4340.2
$ws.Range("K132").Value = 6390.056399999999
$ws.Range("L132").Value = 13020.6
$ws.Range("M132").Value = -3860.056399999999
$ws.Range("N132").Value = -18080.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2919.087
$ws.Range("I7").Value = 1861.125
$ws.Range("J7").Value = 3483.3333
$ws.Range("K7").Value = 1861.125
$ws.Range("L7").Value = 3483.3333
$ws.Range("M7").Value = -1749.125
$ws.Range("N7").Value = -3707.3333
$ws.Range("H16").Value = 4001210
$ws.Range("I16").Value = 4546775
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 4546775
$ws.Range("L16").Value = 400
$ws.Range("M16").Value = -4546605
$ws.Range("N16").Value = -740
$ws.Range("H40").Value = 3154.037
$ws.Range("I40").Value = 1210.4
$ws.Range("K40").Value = 1210.4
$ws.Range("M40").Value = -1074.4
$ws.Range("H46").Value = 2258.4614
$ws.Range("I46").Value = 1800
$ws.Range("J46").Value = 2462.2222
$ws.Range("K46").Value = 1800
$ws.Range("L46").Value = 2462.2222
$ws.Range("M46").Value = -1612
$ws.Range("N46").Value = -2838.2222
$ws.Range("H126").Value = 2919.087
$ws.Range("I126").Value = 1861.125
$ws.Range("J126").Value = 3483.3333
$ws.Range("K126").Value = 5583.375
$ws.Range("L126").Value = 10449.9999
$ws.Range("M126").Value = -3113.375
$ws.Range("N126").Value = -15389.9999
$ws.Range("H132").Value = 5142.298
$ws.Range("I132").Value = 5282.143
$ws.Range("J132").Value = 4734.4165
$ws.Range("K132").Value = 15846.429
$ws.Range("L132").Value = 14203.2495
$ws.Range("M132").Value = -13316.429
$ws.Range("N132").Value = -19263.2495
$ws.Range("H136").Value = 3256.9167
$ws.Range("I136").Value = 1948.5103
$ws.Range("J136").Value = 9085.272000000001
$ws.Range("K136").Value = 5845.5309
$ws.Range("L136").Value = 27255.816
$ws.Range("M136").Value = -3295.5309
$ws.Range("N136").Value = -32355.816

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 22605.264
$ws.Range("J14").Value = 21083.334
$ws.Range("L14").Value = 21083.334
$ws.Range("N14").Value = -21419.334
$ws.Range("H26").Value = 11007
$ws.Range("J26").Value = 11007
$ws.Range("L26").Value = 11007
$ws.Range("N26").Value = -11593
$ws.Range("H81").Value = 1820712.9
$ws.Range("I81").Value = 2502067.8
$ws.Range("J81").Value = 3766.6667
$ws.Range("K81").Value = 5004135.6
$ws.Range("L81").Value = 7533.3334
$ws.Range("M81").Value = -5003074.6
$ws.Range("N81").Value = -9655.3334
$ws.Range("H84").Value = 1820712.9
$ws.Range("I84").Value = 2502067.8
$ws.Range("J84").Value = 3766.6667
$ws.Range("K84").Value = 25020678
$ws.Range("L84").Value = 37666.667
$ws.Range("M84").Value = -25015374
$ws.Range("N84").Value = -48274.667
$ws.Range("H123").Value = 23162.121
$ws.Range("J123").Value = 23162.121
$ws.Range("L123").Value = 23162.121
$ws.Range("N123").Value = -32962.121
$ws.Range("H132").Value = 7938354.5
$ws.Range("I132").Value = 11365445
$ws.Range("J132").Value = 1934.6842
$ws.Range("K132").Value = 34096335
$ws.Range("L132").Value = 5804.0526
$ws.Range("M132").Value = -34093805
$ws.Range("N132").Value = -10864.0526
$ws.Range("H136").Value = 13433.405
$ws.Range("I136").Value = 14771.535
$ws.Range("J136").Value = 1557.5
$ws.Range("K136").Value = 44314.605
$ws.Range("L136").Value = 4672.5
$ws.Range("M136").Value = -41764.605
$ws.Range("N136").Value = -9772.5
